$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EXECUTE column (B) changes
$ws.Range("B5").Value = "NO"
$ws.Range("B10").Value = "YES"
$ws.Range("B11").Value = "YES"

# PARALLEL column (H) changes - standardize case / re-flag each test case
$ws.Range("H2").Value = "YES"
$ws.Range("H3").Value = "NO"
$ws.Range("H4").Value = "YES"
$ws.Range("H5").Value = "YES"
$ws.Range("H6").Value = "NO"
$ws.Range("H7").Value = "NO"
$ws.Range("H8").Value = "NO"
$ws.Range("H9").Value = "NO"
$ws.Range("H10").Value = "YES"
$ws.Range("H11").Value = "YES"
$ws.Range("H12").Value = "YES"
$ws.Range("H13").Value = "YES"
$ws.Range("H14").Value = "YES"

# B14 picks up the same vertical-center alignment as the other EXECUTE cells
$ws.Range("B14").VerticalAlignment = -4108

# Reset the view: scroll back to column A and select B5
$ws.Range("A1").Select()
$ws.Range("B5").Select()
